$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tickers = @(
    "AMZN",
    "NVDA",
    "AMD",
    "LRCX",
    "NFLX",
    "QCOM",
    "AAPL",
    "AMGN",
    "ISRG",
    "GOOG",
    "TXN",
    "INTC",
    "AMAT",
    "CTSH",
    "BIDU",
    "GILD",
    "ADBE",
    "XLNX",
    "QQQ",
    "MSFT",
    "ATVI",
    "CELG",
    "VRTX",
    "TSLA",
    "CSCO",
    "ADSK",
    "CMCSA",
    "BIB",
    "COST",
    "CDNS"
)

for ($i = 0; $i -lt $tickers.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $tickers[$i]
}

$ws.Range("A30").Select()
